$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 12) down to the
# new row (row 13) so the new row's "year" cell picks up the same bold /
# bordered style used by the other year cells in column A.
$ws.Range("A12:O12").Copy()
$ws.Range("A13:O13").PasteSpecial(-4122)

# Fill in the new 2021 data row.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 44601
$ws.Range("C13").Value = 37588
$ws.Range("D13").Value = 36765.66
$ws.Range("E13").Value = 16529.47
$ws.Range("F13").Value = 34335
$ws.Range("G13").Value = 55.1108
$ws.Range("H13").Value = 66958
$ws.Range("I13").Value = 2791
$ws.Range("J13").Value = 47333.44
$ws.Range("K13").Value = 100772
$ws.Range("L13").Value = 4417.6388
$ws.Range("M13").Value = 60356
$ws.Range("N13").Value = 146782
$ws.Range("O13").Value = 89.3484
